$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

function Copy-CellFormat {
    param($srcAddr, $dstAddr)
    $ws.Range($srcAddr).Copy() | Out-Null
    $ws.Range($dstAddr).PasteSpecial(-4122) | Out-Null
}

$cols = @("B","C","D","E","F","H","J","L")

$newRows = @(
  @{ Row=265; B="riesgos"; C=41807.97292824074; D=41807.97292824074; E="riesgos salio del sistema"; F=0.0; H=0.0; J=0.0; L=" " },
  @{ Row=266; B="riesgos"; C=41808.01106481482; D=41808.01106481482; E="riesgos ingreso al sistema"; F=0.0; H=0.0; J=0.0; L=" " },
  @{ Row=267; B="riesgos"; C=41808.01798611111; D=41808.01798611111; E="riesgos ingreso al sistema"; F=0.0; H=0.0; J=0.0; L=" " },
  @{ Row=268; B="riesgos"; C=41808.01913194444; D=41808.01913194444; E="riesgos ingreso al sistema"; F=0.0; H=0.0; J=0.0; L=" " },
  @{ Row=269; B="riesgos"; C=41808.01954861111; D=41808.01954861111; E="Genero reporte comparativo de ejemplo y ejemplo 2"; F=0.0; H=0.0; J=0.0; L=" " },
  @{ Row=270; B="riesgos"; C=41808.0262037037; D=41808.0262037037; E="riesgos ingreso al sistema"; F=0.0; H=0.0; J=0.0; L=" " },
  @{ Row=271; B="riesgos"; C=41808.026608796295; D=41808.026608796295; E="Genero reporte comparativo de ejemplo y ejemplo 2"; F=0.0; H=0.0; J=0.0; L=" " },
  @{ Row=272; B="riesgos"; C=41808.05274305555; D=41808.05274305555; E="riesgos ingreso al sistema"; F=0.0; H=0.0; J=0.0; L=" " },
  @{ Row=273; B="riesgos"; C=41808.053148148145; D=41808.053148148145; E="Genero reporte comparativo de ejemplo y ejemplo 2"; F=0.0; H=0.0; J=0.0; L=" " },
  @{ Row=274; B="riesgos"; C=41809.98532407408; D=41809.98532407408; E="riesgos ingreso al sistema"; F=0.0; H=0.0; J=0.0; L=" " },
  @{ Row=275; B="riesgos"; C=41809.98537037037; D=41809.98537037037; E="riesgos salio del sistema"; F=0.0; H=0.0; J=0.0; L=" " },
  @{ Row=276; B="riesgos"; C=41809.987222222226; D=41809.987222222226; E="riesgos ingreso al sistema"; F=0.0; H=0.0; J=0.0; L=" " },
  @{ Row=277; B="riesgos"; C=41809.987604166665; D=41809.987604166665; E="riesgos salio del sistema"; F=0.0; H=0.0; J=0.0; L=" " },
  @{ Row=278; B="riesgos"; C=41809.98829861111; D=41809.98829861111; E="riesgos ingreso al sistema"; F=0.0; H=0.0; J=0.0; L=" " },
  @{ Row=279; B="riesgos"; C=41809.98836805556; D=41809.98836805556; E="riesgos salio del sistema"; F=0.0; H=0.0; J=0.0; L=" " },
  @{ Row=280; B="riesgos"; C=41809.99302083333; D=41809.99302083333; E="riesgos ingreso al sistema"; F=0.0; H=0.0; J=0.0; L=" " },
  @{ Row=281; B="riesgos"; C=41809.9930787037; D=41809.9930787037; E="riesgos salio del sistema"; F=0.0; H=0.0; J=0.0; L=" " },
  @{ Row=282; B="riesgos"; C=41809.99353009259; D=41809.99353009259; E="riesgos ingreso al sistema"; F=0.0; H=0.0; J=0.0; L=" " },
  @{ Row=283; B="riesgos"; C=41809.99375; D=41809.99375; E="riesgos salio del sistema"; F=0.0; H=0.0; J=0.0; L=" " },
  @{ Row=284; B="riesgos"; C=41809.993842592594; D=41809.993842592594; E="riesgos ingreso al sistema"; F=0.0; H=0.0; J=0.0; L=" " },
  @{ Row=285; B="riesgos"; C=41809.99387731482; D=41809.99387731482; E="riesgos salio del sistema"; F=0.0; H=0.0; J=0.0; L=" " },
  @{ Row=286; B="riesgos"; C=41809.99619212963; D=41809.99619212963; E="riesgos ingreso al sistema"; F=0.0; H=0.0; J=0.0; L=" " },
  @{ Row=287; B="riesgos"; C=41809.99627314815; D=41809.99627314815; E="riesgos salio del sistema"; F=0.0; H=0.0; J=0.0; L=" " },
  @{ Row=288; B="riesgos"; C=41809.996516203704; D=41809.996516203704; E="riesgos ingreso al sistema"; F=0.0; H=0.0; J=0.0; L=" " },
  @{ Row=289; B="riesgos"; C=41809.99662037037; D=41809.99662037037; E="riesgos salio del sistema"; F=0.0; H=0.0; J=0.0; L=" " },
  @{ Row=290; B="riesgos"; C=41809.99873842593; D=41809.99873842593; E="riesgos ingreso al sistema"; F=0.0; H=0.0; J=0.0; L=" " },
  @{ Row=291; B="riesgos"; C=41809.99878472222; D=41809.99878472222; E="riesgos salio del sistema"; F=0.0; H=0.0; J=0.0; L=" " },
  @{ Row=292; B="riesgos"; C=41809.99900462963; D=41809.99900462963; E="riesgos ingreso al sistema"; F=0.0; H=0.0; J=0.0; L=" " }
)

foreach ($rd in $newRows) {
    $r = $rd.Row
    foreach ($col in $cols) {
        Copy-CellFormat "${col}264" "${col}${r}"
    }
    $ws.Range("B$r").Value2 = $rd.B
    $ws.Range("C$r").Value2 = $rd.C
    $ws.Range("D$r").Value2 = $rd.D
    $ws.Range("E$r").Value2 = $rd.E
    $ws.Range("F$r").Value2 = $rd.F
    $ws.Range("H$r").Value2 = $rd.H
    $ws.Range("J$r").Value2 = $rd.J
    $ws.Range("L$r").Value2 = $rd.L
}

$excel.CutCopyMode = $false
Write-Host "Done adding rows 265-292"